$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("personas")
$lo = $ws.ListObjects.Item("Table1")

# Rename the SharePoint-prefixed skill/school columns (set header cell text
# directly since the table syncs its column name from the header row).
$ws.Range("L1").Value = "Skills"
$ws.Range("M1").Value = "School"

# The DisplayName column should now be computed from Title/FirstName/LastName
# instead of being a manually-typed value.
$dispCol = $lo.ListColumns.Item("DisplayName")
$dispRng = $dispCol.DataBodyRange
$dispRng.Formula = "=_xlfn.CONCAT(Table1[[#This Row],[Title]], "" "", Table1[[#This Row],[FirstName]], "" "", Table1[[#This Row],[LastName]])"

# Flag the (now calculated) column visually with the "Check Cell" style.
$dispRng.Style = "Check Cell"

# Restore the view so the newly important columns are visible/selected.
$ws.Application.ActiveWindow.ScrollColumn = 7
$sel = $ws.Range("M1")
$sel.Select()
